$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''25.663.14'
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").Value = '''1.691.12'
$ws.Range("E3").Value = '  -2.76%  '
$ws.Range("D4").Value = '''1.006'
$ws.Range("E4").Value = '  +0.51%  '
$ws.Range("D5").Value = '''241.22'
$ws.Range("E5").Value = '  +1.54%  '
$ws.Range("D6").Value = '''1.005'
$ws.Range("E6").Value = '  +0.44%  '
$ws.Range("D7").Value = '''0.4871'
$ws.Range("E7").Value = '  -5.82%  '
$ws.Range("D8").Value = '''0.2655'
$ws.Range("E8").Value = '  -3.36%  '
$ws.Range("D9").Value = '''0.06047'
$ws.Range("E9").Value = '  -1.77%  '
$ws.Range("D10").Value = '''1.716.26'
$ws.Range("E10").Value = '  -1.34%  '
$ws.Range("D11").Value = '''0.07171'
$ws.Range("E11").Value = '  -0.20%  '
$ws.Range("D12").Value = '''0.6327'
$ws.Range("E12").Value = '  -1.36%  '
$ws.Range("D13").Value = '''14.68'
$ws.Range("E13").Value = '  -1.58%  '
$ws.Range("D14").Value = '''4.650'
$ws.Range("E14").Value = '  +1.13%  '
$ws.Range("D15").Value = '''74.68'
$ws.Range("E15").Value = '  -3.54%  '
$ws.Range("B16").Value = 'Dai'
$ws.Range("C16").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D16").Value = '''1.005'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").Value = '''1.006'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").Value = '''25.670.89'
$ws.Range("E18").Value = '  -0.83%  '
$ws.Range("D19").Value = '''11.58'
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("D20").Value = '''0.000006687'
$ws.Range("E20").Value = '  -1.15%  '
$ws.Range("D21").Value = '''1.927.06'
$ws.Range("E21").Value = '  -1.99%  '
$ws.Range("D22").Value = '''4.485'
$ws.Range("E22").Value = '  +4.86%  '
$ws.Range("D23").Value = '''8.666'
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("D24").Value = '''5.333'
$ws.Range("E24").Value = '  +1.38%  '
$ws.Range("D25").Value = '''133.63'
$ws.Range("E25").Value = '  -3.87%  '
$ws.Range("D26").Value = '''14.96'
$ws.Range("E26").Value = '  -1.18%  '
$ws.Range("D27").Value = '''1.396'
$ws.Range("E27").Value = '  -7.77%  '
$ws.Range("D28").Value = '''1.736'
$ws.Range("E28").Value = '  -1.08%  '
$ws.Range("D29").Value = '''103.40'
$ws.Range("E29").Value = '  -2.09%  '
$ws.Range("D30").Value = '''3.860'
$ws.Range("E30").Value = '  -1.21%  '
$ws.Range("D31").Value = '''0.08013'
$ws.Range("E31").Value = '  -3.07%  '
$ws.Range("D32").Value = '''3.578'
$ws.Range("E32").Value = '  -3.05%  '
$ws.Range("D33").Value = '''0.04634'
$ws.Range("E33").Value = '  +0.79%  '
$ws.Range("D34").Value = '''2.672'
$ws.Range("E34").Value = '  +1.05%  '
$ws.Range("D35").Value = '''0.9678'
$ws.Range("E35").Value = '  -1.96%  '
$ws.Range("D36").Value = '''0.5909'
$ws.Range("E36").Value = '  -4.16%  '
$ws.Range("D37").Value = '''2.684'
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").Value = '''0.01571'
$ws.Range("E38").Value = '  -1.92%  '
$ws.Range("D39").Value = '''0.8466'
$ws.Range("E39").Value = '  +14.56%  '
$ws.Range("E40").Value = '  +0.53%  '
$ws.Range("D41").Value = '''1.891'
$ws.Range("E41").Value = '  -1.59%  '
$ws.Range("D42").Value = '''99.61'
$ws.Range("E42").Value = '  +2.01%  '
$ws.Range("D43").Value = '''0.3775'
$ws.Range("E43").Value = '  -1.49%  '
$ws.Range("D44").Value = '''4.928'
$ws.Range("E44").Value = '  -1.04%  '
$ws.Range("D45").Value = '''0.1156'
$ws.Range("E45").Value = '  +2.41%  '
$ws.Range("D46").Value = '''6.147'
$ws.Range("E46").Value = '  -1.12%  '
$ws.Range("D47").Value = '''0.05226'
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("D48").Value = '''54.35'
$ws.Range("E48").Value = '  -0.64%  '
$ws.Range("D49").Value = '''29.92'
$ws.Range("E49").Value = '  -1.54%  '
$ws.Range("D50").Value = '''7.494'
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("D51").Value = '''0.3383'
$ws.Range("E51").Value = '  -0.48%  '

# Re-apply the original (unstyled/general) look to the data cells so that the
# quote-prefix formatting Excel applied above does not linger as a style change.
$ws.Range("D2:E51").Style = $ws.Range("B2").Style
